$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the French intro text in E2: remove the extra blank line between
# the two sentences (double newline -> single newline)
$ws.Range("E2").Value = "Vous allez maintenant faire les blocs de test`nIl y a 4 blocs, avec une pause entre les deux"

# Rename the header labels in row 1 (uppercase the language suffix and
# rename the "test_msg" labels to "block_msg")
$ws.Range("B1").Value = "block_pause_msg_EN"
$ws.Range("D1").Value = "block_pause_msg_ES"
$ws.Range("F1").Value = "block_pause_msg_FR"
$ws.Range("A1").Value = "block_msg_EN"
$ws.Range("C1").Value = "block_msg_ES"
$ws.Range("E1").Value = "block_msg_FR"
